# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets,
# mirroring the data refresh recorded in the commit diff.

$wb = $excel.ActiveWorkbook

# row => new F-column value, for worksheet "展览"
$exhibitionUpdates = @{
    2  = 1681
    3  = 9124
    7  = 1376
    9  = 60
    10 = 97
    11 = 5915
    15 = 4463
    18 = 1151
    19 = 30
    20 = 340
    21 = 27
    22 = 260
    24 = 2763
}

# row => new F-column value, for worksheet "全部类型"
$allTypesUpdates = @{
    2  = 1681
    3  = 9124
    8  = 1376
    10 = 60
    11 = 97
    12 = 5915
    16 = 4463
    19 = 1151
    20 = 30
    21 = 340
    22 = 27
    23 = 260
    25 = 2763
}

$wsExhibition = $wb.Worksheets.Item("展览")
foreach ($row in $exhibitionUpdates.Keys) {
    $wsExhibition.Range("F$row").Value = $exhibitionUpdates[$row]
}

$wsAllTypes = $wb.Worksheets.Item("全部类型")
foreach ($row in $allTypesUpdates.Keys) {
    $wsAllTypes.Range("F$row").Value = $allTypesUpdates[$row]
}
